$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (48) used the "final row" date format (YYYY-MM-DD).
# Since we're appending a new last row, row 48 reverts to the regular
# date/time format used by all the other data rows.
$ws.Range("A48").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 49: next day's data (daily update).
$ws.Range("A49").Value = 45634
$ws.Range("A49").NumberFormat = "YYYY-MM-DD"
$ws.Range("B49").Value = 127
$ws.Range("C49").Value = 113
$ws.Range("D49").Value = 118
